# 🔄 MAJ automatique BRVM via GitHub Actions
# Refresh the "Recommandations" and "Top_YTD" sheets with the latest BRVM
# market data (rankings re-sorted by "Variation Totale (%)", a new line
# added for SAFCA CI (SAFC), and the Top_YTD progression figures updated).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Recommandations"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Recommandations")

$reco = @(
    @(2,  "BRVM - SERVICES PUBLICS",                 0, 6, 2445.53, 100.18,              "🟡 Observer", "➖ Neutre"),
    @(3,  "BRVM - AUTRES SECTEURS",                   0, 3, 1905.63, 632.35,              "🟡 Observer", "➖ Neutre"),
    @(4,  "UNIWAX CI",                                0, 3, 1865,    665,                 "🟡 Observer", "➖ Neutre"),
    @(5,  "NEI-CEDA CI",                              0, 3, 1750,    600,                 "🟡 Observer", "➖ Neutre"),
    @(6,  "SETAO CI",                                 0, 3, 1735,    595,                 "🟡 Observer", "➖ Neutre"),
    @(7,  "AIR LIQUIDE CI",                           0, 3, 1635,    545,                 "🟡 Observer", "➖ Neutre"),
    @(8,  "CFAO MOTORS CI",                           0, 2, 1315,    650,                 "🟡 Observer", "➖ Neutre"),
    @(9,  "BRVM - DISTRIBUTION",                      0, 3, 1116.73, 372.68,              "🟡 Observer", "➖ Neutre"),
    @(10, "BRVM - TRANSPORT",                         0, 3, 1088.02, 367.16,              "🟡 Observer", "➖ Neutre"),
    @(11, "SAFCA CI",                                 0, 1, 995,     995,                 "🟡 Observer", "➖ Neutre"),
    @(12, "BRVM - AGRICULTURE",                       0, 3, 973.25,  323.81,              "🟡 Observer", "➖ Neutre"),
    @(13, "BRVM - INDUSTRIE",                         0, 3, 793.29,  267.92,              "🟡 Observer", "➖ Neutre"),
    @(14, "BRVM - CONSOMMATION DE BASE",              0, 3, 655.01,  220.49,              "🟡 Observer", "➖ Neutre"),
    @(15, "BRVM-PRINCIPAL",                           0, 3, 572.16,  191.99,              "🟡 Observer", "➖ Neutre"),
    @(16, "BRVM - INDUSTRIELS",                       0, 3, 416.35,  140.04,              "🟡 Observer", "➖ Neutre"),
    @(17, "BRVM-PRESTIGE",                            0, 3, 391.96,  131.38,              "🟡 Observer", "➖ Neutre"),
    @(18, "BRVM - FINANCES",                          0, 3, 369.1,   123.72,              "🟡 Observer", "➖ Neutre"),
    @(19, "BRVM - SERVICES FINANCIERS",               0, 3, 362.74,  121.59,              "🟡 Observer", "➖ Neutre"),
    @(20, "BRVM - ENERGIE",                           0, 3, 331.96,  110.44,              "🟡 Observer", "➖ Neutre"),
    @(21, "BRVM - CONSOMMATION DISCRETIONNAIRE",      0, 3, 320.27,  107.38,              "🟡 Observer", "➖ Neutre"),
    @(22, "BRVM - TELECOMMUNICATIONS",                0, 3, 280.22,  93.70999999999999,   "🟡 Observer", "➖ Neutre"),
    @(23, "BERNABE CI (BNBC)",                        2, 0, 14.28,   7.21,                "🟡 Observer", "➖ Neutre"),
    @(24, "UNIWAX CI (UNXC)",                         2, 0, 14.16,   7.26,                "🟡 Observer", "➖ Neutre"),
    @(25, "SUCRIVOIRE (SCRC)",                        1, 0, 6.76,    6.76,                "🟡 Observer", "➖ Neutre"),
    @(26, "AFRICA GLOBAL LOGISTICS CI (SDSC)",        1, 0, 3.81,    3.81,                "🟡 Observer", "➖ Neutre"),
    @(27, "ECOBANK COTE D''IVOIRE (ECOC)",            1, 1, 3.12,    -0.42,               "🟡 Observer", "👀 À surveiller"),
    @(28, "CFAO MOTORS CI (CFAC)",                    1, 0, 3.05,    3.05,                "🟡 Observer", "➖ Neutre"),
    @(29, "ORANGE COTE D'IVOIRE (ORAC)",              1, 0, 2.93,    2.93,                "🟡 Observer", "➖ Neutre"),
    @(30, "SETAO CI (STAC)",                          1, 1, 2.8,     6.25,                "🟡 Observer", "👀 À surveiller"),
    @(31, "ONATEL BF (ONTBF)",                        1, 0, 1.96,    1.96,                "🟡 Observer", "➖ Neutre"),
    @(32, "ORAGROUP TOGO (ORGT)",                     1, 1, 1.89,    7.26,                "🟡 Observer", "👀 À surveiller"),
    @(33, "SICABLE CI (CABC)",                        1, 1, 0.44,    -2.25,               "🟡 Observer", "👀 À surveiller"),
    @(34, "ECOBANK TRANS. INCORP. TG (ETIT)",         1, 1, 0.37,    6.25,                "🟡 Observer", "👀 À surveiller"),
    @(36, "TOTALENERGIES MARKETING SN (TTLS)",        0, 1, -0.2,    -0.2,                "🟡 Observer", "➖ Neutre"),
    @(37, "LOTERIE NATIONALE DU BENIN (LNBB)",        0, 1, -0.55,   -0.55,               "🟡 Observer", "➖ Neutre"),
    @(38, "BICI CI (BICC)",                           0, 1, -1.81,   -1.81,               "🟡 Observer", "➖ Neutre"),
    @(39, "SOLIBRA CI (SLBC)",                        0, 1, -1.97,   -1.97,               "🟡 Observer", "➖ Neutre"),
    @(40, "SOCIETE GENERALE COTE D'IVOIRE (SGBC)",    0, 1, -2.02,   -2.02,               "🟡 Observer", "➖ Neutre"),
    @(41, "VIVO ENERGY CI (SHEC)",                    0, 1, -2.53,   -2.53,               "🟡 Observer", "➖ Neutre"),
    @(42, "SOGB CI (SOGC)",                           0, 1, -3.16,   -3.16,               "🟡 Observer", "➖ Neutre"),
    @(43, "SERVAIR ABIDJAN CI (ABJC)",                0, 1, -3.68,   -3.68,               "🟡 Observer", "➖ Neutre"),
    @(44, "SAFCA CI (SAFC)",                          1, 2, -7.73,   7.04,                "🟡 Observer", "👀 À surveiller")
)

foreach ($line in $reco) {
    $r = $line[0]
    $ws1.Cells.Item($r, 1).Value = $line[1]
    $ws1.Cells.Item($r, 2).Value = $line[2]
    $ws1.Cells.Item($r, 3).Value = $line[3]
    $ws1.Cells.Item($r, 4).Value = $line[4]
    $ws1.Cells.Item($r, 5).Value = $line[5]
    $ws1.Cells.Item($r, 6).Value = $line[6]
    $ws1.Cells.Item($r, 7).Value = $line[7]
}

# ---------------------------------------------------------------------
# Sheet 2: "Top_YTD"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Top_YTD")

$topYtd = @(
    @(2,  "BRVM - SERVICES PUBLICS", 433024.02),
    @(3,  "BRVM - AUTRES SECTEURS",  39640.05),
    @(4,  "UNIWAX CI",               37354.4),
    @(5,  "NEI-CEDA CI",             31765.75),
    @(6,  "SETAO CI",                31091.6),
    @(7,  "AIR LIQUIDE CI",          26733.61),
    @(8,  "BRVM - DISTRIBUTION",     10431.67),
    @(9,  "BRVM - TRANSPORT",        9801.48),
    @(10, "BRVM - AGRICULTURE",      7544.81),
    @(11, "CFAO MOTORS CI",          5637.5)
)

foreach ($line in $topYtd) {
    $r = $line[0]
    $ws2.Cells.Item($r, 1).Value = $line[1]
    $ws2.Cells.Item($r, 2).Value = $line[2]
}
